$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Swap the "(3-4)" row (old row 4) with the "(5-6) Centre" row (old row 3),
# so that "(3-4)" appears above "(5-6) Centre" in the table.
$row3vals = $ws.Range("A3:L3").Value2
$row4vals = $ws.Range("A4:L4").Value2
$ws.Range("A3:L3").Value2 = $row4vals
$ws.Range("A4:L4").Value2 = $row3vals

# Re-apply consistent formatting (matching row 2's style) across rows 3-6.
$null = $ws.Range("A2:L2").Copy()
$null = $ws.Range("A3:L6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The custom row height follows the "(5-6) Centre" data, which is now row 4;
# row 3 goes back to the default row height.
$ws.Rows.Item(4).RowHeight = 23.25
$null = $ws.Rows.Item(3).AutoFit()

# Add a new (mostly empty) row 7 below the table, matching the formatting
# used by the rest of the table in column B.
$null = $ws.Range("B3").Copy()
$null = $ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Activate()
$null = $ws.Range("M15").Select()
